# Apply updated cryptocurrency price/volume snapshot (and restore the
# PancakeSwap/Decentraland row ordering) to Sheet1, per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value reads as a plain number (e.g. "313.50") need a
# leading apostrophe so Excel keeps storing them as text, matching the
# workbook's existing inline-string cells instead of silently coercing
# them (and dropping trailing zeros / switching to scientific notation).
$updates = @(
    @{ Cell = 'D2'; Value = '27.277.49' }
    @{ Cell = 'E2'; Value = '  +0.58%  ' }
    @{ Cell = 'D3'; Value = '1.771.94' }
    @{ Cell = 'E3'; Value = '  +3.54%  ' }
    @{ Cell = 'E4'; Value = '  +0.16%  ' }
    @{ Cell = 'D5'; Value = '''313.50' }
    @{ Cell = 'E5'; Value = '  +2.07%  ' }
    @{ Cell = 'D6'; Value = '''1.001' }
    @{ Cell = 'E6'; Value = '  +0.06%  ' }
    @{ Cell = 'D7'; Value = '''0.5270' }
    @{ Cell = 'E7'; Value = '  +10.13%  ' }
    @{ Cell = 'D8'; Value = '''0.3625' }
    @{ Cell = 'E8'; Value = '  +5.34%  ' }
    @{ Cell = 'D9'; Value = '''42.62' }
    @{ Cell = 'E9'; Value = '  +1.68%  ' }
    @{ Cell = 'D10'; Value = '''0.07358' }
    @{ Cell = 'E10'; Value = '  +0.93%  ' }
    @{ Cell = 'D11'; Value = '''1.091' }
    @{ Cell = 'E11'; Value = '  +4.04%  ' }
    @{ Cell = 'D12'; Value = '''1.001' }
    @{ Cell = 'E12'; Value = '  +0.08%  ' }
    @{ Cell = 'D13'; Value = '''20.49' }
    @{ Cell = 'E13'; Value = '  +3.20%  ' }
    @{ Cell = 'D14'; Value = '''6.065' }
    @{ Cell = 'E14'; Value = '  +3.60%  ' }
    @{ Cell = 'D15'; Value = '1.771.92' }
    @{ Cell = 'E15'; Value = '  +3.97%  ' }
    @{ Cell = 'D16'; Value = '''6.955' }
    @{ Cell = 'E16'; Value = '  +1.66%  ' }
    @{ Cell = 'D17'; Value = '''88.45' }
    @{ Cell = 'E17'; Value = '  -0.78%  ' }
    @{ Cell = 'D18'; Value = '''0.00001044' }
    @{ Cell = 'E18'; Value = '  +0.38%  ' }
    @{ Cell = 'D19'; Value = '''0.06421' }
    @{ Cell = 'E19'; Value = '  +1.13%  ' }
    @{ Cell = 'E20'; Value = '  +0.10%  ' }
    @{ Cell = 'D21'; Value = '''16.73' }
    @{ Cell = 'E21'; Value = '  +1.62%  ' }
    @{ Cell = 'D22'; Value = '''5.826' }
    @{ Cell = 'E22'; Value = '  +3.98%  ' }
    @{ Cell = 'D23'; Value = '27.373.22' }
    @{ Cell = 'E23'; Value = '  +0.84%  ' }
    @{ Cell = 'D24'; Value = '''11.30' }
    @{ Cell = 'E24'; Value = '  +4.16%  ' }
    @{ Cell = 'D25'; Value = '''2.083' }
    @{ Cell = 'E25'; Value = '  -0.91%  ' }
    @{ Cell = 'D26'; Value = '''154.10' }
    @{ Cell = 'E26'; Value = '  -0.89%  ' }
    @{ Cell = 'D27'; Value = '''20.11' }
    @{ Cell = 'E27'; Value = '  +2.53%  ' }
    @{ Cell = 'D28'; Value = '''2.343' }
    @{ Cell = 'E28'; Value = '  +12.12%  ' }
    @{ Cell = 'D29'; Value = '1.974.69' }
    @{ Cell = 'E29'; Value = '  +5.32%  ' }
    @{ Cell = 'D30'; Value = '''121.13' }
    @{ Cell = 'E30'; Value = '  +1.40%  ' }
    @{ Cell = 'D31'; Value = '''1.059' }
    @{ Cell = 'E31'; Value = '  +4.48%  ' }
    @{ Cell = 'D32'; Value = '''0.09804' }
    @{ Cell = 'E32'; Value = '  +6.59%  ' }
    @{ Cell = 'D33'; Value = '''5.534' }
    @{ Cell = 'E33'; Value = '  +4.14%  ' }
    @{ Cell = 'D34'; Value = '''3.620' }
    @{ Cell = 'E34'; Value = '  +0.97%  ' }
    @{ Cell = 'D35'; Value = '''0.02225' }
    @{ Cell = 'E35'; Value = '  +1.32%  ' }
    @{ Cell = 'D36'; Value = '''0.05959' }
    @{ Cell = 'E36'; Value = '  +2.24%  ' }
    @{ Cell = 'D37'; Value = '''11.15' }
    @{ Cell = 'E37'; Value = '  +0.74%  ' }
    @{ Cell = 'D38'; Value = '''4.836' }
    @{ Cell = 'E38'; Value = '  +2.05%  ' }
    @{ Cell = 'D39'; Value = '''0.6131' }
    @{ Cell = 'E39'; Value = '  +4.44%  ' }
    @{ Cell = 'D40'; Value = '''0.2020' }
    @{ Cell = 'E40'; Value = '  +1.36%  ' }
    @{ Cell = 'D41'; Value = '''1.435' }
    @{ Cell = 'E41'; Value = '  +2.50%  ' }
    @{ Cell = 'D42'; Value = '''8.044' }
    @{ Cell = 'E42'; Value = '  +7.71%  ' }
    @{ Cell = 'D43'; Value = '''1.145' }
    @{ Cell = 'E43'; Value = '  +3.33%  ' }
    @{ Cell = 'D44'; Value = '''13.07' }
    @{ Cell = 'E44'; Value = '  +3.20%  ' }
    @{ Cell = 'B45'; Value = 'PancakeSwap' }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake' }
    @{ Cell = 'D45'; Value = '''3.631' }
    @{ Cell = 'E45'; Value = '  +2.17%  ' }
    @{ Cell = 'B46'; Value = 'Decentraland' }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana' }
    @{ Cell = 'D46'; Value = '''0.5751' }
    @{ Cell = 'E46'; Value = '  +2.22%  ' }
    @{ Cell = 'D47'; Value = '''120.97' }
    @{ Cell = 'E47'; Value = '  +2.80%  ' }
    @{ Cell = 'D48'; Value = '''1.884' }
    @{ Cell = 'E48'; Value = '  +2.27%  ' }
    @{ Cell = 'D49'; Value = '''1.109' }
    @{ Cell = 'E49'; Value = '  +2.25%  ' }
    @{ Cell = 'D50'; Value = '''0.06709' }
    @{ Cell = 'E50'; Value = '  +1.22%  ' }
    @{ Cell = 'D51'; Value = '''70.58' }
    @{ Cell = 'E51'; Value = '  +1.21%  ' }
)

foreach ($update in $updates) {
    $ws.Range($update.Cell).Value = $update.Value
}
